# Update gh-pages output data: increment "想去人数" (want-to-go count) values
# for several events on the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5621
$ws1.Range("F6").Value = 932
$ws1.Range("F8").Value = 2517
$ws1.Range("F10").Value = 144
$ws1.Range("F11").Value = 9
$ws1.Range("F13").Value = 18
$ws1.Range("F14").Value = 2361
$ws1.Range("F15").Value = 345

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5621
$ws4.Range("F8").Value = 932
$ws4.Range("F10").Value = 2517
$ws4.Range("F12").Value = 144
$ws4.Range("F13").Value = 9
$ws4.Range("F16").Value = 18
$ws4.Range("F17").Value = 2361
$ws4.Range("F18").Value = 345
